$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated loading_percent values (rows 2-25; columns B,C,D,F,G,H,I,J)
# Each row entry: row number, then values for columns B,C,D,F,G,H,I,J
$rowsData = @(
    ,@(2, 16.57652354381622, 12.0179391666997, 10.64930469477278, 32.87581516399901, 33.54123221659911, 15.24229421831632, 20.65766956351989, 11.33742494686702)
    ,@(3, 15.87523957075624, 11.38083011061292, 10.5698421961238, 32.87907548986671, 33.51339753248437, 15.30743810019331, 20.82056591486488, 11.32127608191666)
    ,@(4, 15.42991651138738, 10.97173728011472, 10.52223316956167, 32.89377933704076, 33.51443709199092, 15.35177067930391, 20.92746547961286, 11.31377666514362)
    ,@(5, 15.24498941369897, 10.80069448156167, 10.50314418153018, 32.90295194202901, 33.51939528657008, 15.37092096523281, 20.97275093195296, 11.31132963187937)
    ,@(6, 15.21408172534895, 10.77203672541796, 10.49999374761643, 32.90466670182168, 33.52049161964902, 15.37416619102398, 20.98037440787811, 11.3109601232862)
    ,@(7, 15.42743615202942, 10.96944784038692, 10.52197444571892, 32.89389018360743, 33.51448563719535, 15.35202456334304, 20.92806924804166, 11.31374119597809)
    ,@(8, 16.33791862184389, 11.80207179585775, 10.62167002202078, 32.8742959039737, 33.5278593470112, 15.2638530810053, 20.71240379995712, 11.33135640293135)
    ,@(9, 17.99706543795873, 13.2868334650534, 10.82586410209789, 32.93711035114306, 33.69868289306758, 15.12557604950412, 20.34440777578984, 11.3849852557477)
    ,@(10, 19.12841870697823, 14.28119070131794, 10.98024626536915, 33.04538777463353, 33.91284136030261, 15.04541524143807, 20.10798943351103, 11.43588061152636)
    ,@(11, 19.6223314346237, 14.71166706370346, 11.05122016442852, 33.10815474778799, 34.0294672635288, 15.01366862823244, 20.00791447009023, 11.4614905189759)
    ,@(12, 19.80626313905879, 14.87147597718161, 11.07818564069411, 33.13386178312426, 34.07637600753094, 15.00233093922991, 19.97110275147302, 11.47153755365657)
    ,@(13, 19.76678983485794, 14.83720154444431, 11.07237446097205, 33.12823920111795, 34.06615162511418, 15.00474220352328, 19.97898241898097, 11.469358290221)
    ,@(14, 19.63752639901472, 14.72487913138886, 11.05343694529576, 33.11023086670192, 34.03327158021555, 15.01272212912155, 20.00486415494234, 11.46231012362856)
    ,@(15, 19.55794157682818, 14.6556597150978, 11.04184825499784, 33.09945251071567, 34.01348847281734, 15.01769930987906, 20.02085899898064, 11.45803824991767)
    ,@(16, 19.09571228266009, 14.25261277148404, 10.97562141192457, 33.04155739051497, 33.90560478108567, 15.04758538644264, 20.11468078047489, 11.43425603947006)
    ,@(17, 18.80674535206898, 13.99971079102198, 10.93517142060116, 33.00949912219409, 33.8443309219854, 15.06713210543075, 20.17415833426133, 11.4202929342272)
    ,@(18, 18.63859239271939, 13.85219319237984, 10.91197682533319, 32.99233225443709, 33.81089687270909, 15.07881871022111, 20.20907102386599, 11.41249330660562)
    ,@(19, 18.58132826037478, 13.80189533038557, 10.9041363010508, 32.98673842494139, 33.79988776427993, 15.0828516559031, 20.22101226667478, 11.40989238293662)
    ,@(20, 18.83770878614809, 14.02684580096939, 10.93947015474263, 33.01278013244258, 33.85066649479285, 15.06500534713785, 20.16775404236089, 11.42175539218784)
    ,@(21, 19.6755792581673, 14.75795829151006, 11.0589970746044, 33.1154677909639, 34.04285491494862, 15.01035962316447, 19.99723254810037, 11.46437090209826)
    ,@(22, 20.20505055929059, 15.21709296324635, 11.13762644036836, 33.19387501402013, 34.18444853105425, 14.97863518233315, 19.89211302201067, 11.49425524764789)
    ,@(23, 19.92415495648313, 14.97377077306243, 11.09561957385421, 33.15099643284778, 34.10742186625743, 14.99520027861899, 19.94763508860559, 11.47812094090737)
    ,@(24, 18.82371650651677, 14.01458466446622, 10.93752650679254, 33.01129285062165, 33.84779659582886, 15.06596545656262, 20.17064718615046, 11.42109350461139)
    ,@(25, 17.56294244846251, 12.90174042808348, 10.7697892340432, 32.90922376950051, 33.63690891674104, 15.1592429919309, 20.43803363356785, 11.36844666147552)
)

$cols = @(2, 3, 4, 6, 7, 8, 9, 10)  # B, C, D, F, G, H, I, J

foreach ($entry in $rowsData) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $c = $cols[$i]
        $val = $entry[$i + 1]
        $ws.Cells.Item($r, $c).Value = $val
    }
}
